$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget Roboter")

# Update the price for the "Mechanik" / "Schienen m. Wagen" row (C4) from 41 to 70
$ws.Range("C4").Value = 70

# Recalculate formulas (F4 = E4+C4, I11 = SUM(F4:F9)) so dependent values update
$excel.Calculate()

# Update the active selection to match the new selection in the diff
$ws.Range("E16").Select()

$wb.Save()
